$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 500
$ws.Range("I31").Value = 500
$ws.Range("K31").Value = 1500
$ws.Range("M31").Value = -1270
$ws.Range("H33").Value = 7473.1333
$ws.Range("I33").Value = 322
$ws.Range("J33").Value = 18199.834
$ws.Range("K33").Value = 322
$ws.Range("L33").Value = 18199.834
$ws.Range("M33").Value = -93
$ws.Range("N33").Value = -18657.834
$ws.Range("H74").Value = 6082.5884
$ws.Range("I74").Value = 9800
$ws.Range("J74").Value = 4533.6665
$ws.Range("K74").Value = 9800
$ws.Range("L74").Value = 4533.6665
$ws.Range("M74").Value = -8864
$ws.Range("N74").Value = -6405.6665
$ws.Range("H76").Value = 2852118.2
$ws.Range("I76").Value = 4632442
$ws.Range("J76").Value = 3600
$ws.Range("K76").Value = 4632442
$ws.Range("L76").Value = 3600
$ws.Range("M76").Value = -4632127
$ws.Range("N76").Value = -4230
$ws.Range("H77").Value = 6082.5884
$ws.Range("I77").Value = 9800
$ws.Range("J77").Value = 4533.6665
$ws.Range("K77").Value = 49000
$ws.Range("L77").Value = 22668.3325
$ws.Range("M77").Value = -44320
$ws.Range("N77").Value = -32028.3325
$ws.Range("H79").Value = 2852118.2
$ws.Range("I79").Value = 4632442
$ws.Range("J79").Value = 3600
$ws.Range("K79").Value = 4632442
$ws.Range("L79").Value = 3600
$ws.Range("M79").Value = -4631350
$ws.Range("N79").Value = -5784
$ws.Range("H132").Value = 2526834.2
$ws.Range("I132").Value = 2675447.8
$ws.Range("J132").Value = 406
$ws.Range("K132").Value = 8026343.399999999
$ws.Range("L132").Value = 1218
$ws.Range("M132").Value = -8023813.399999999
$ws.Range("N132").Value = -6278
$ws.Range("H137").Value = 1387.6923
$ws.Range("I137").Value = 1224
$ws.Range("J137").Value = 1933.3334
$ws.Range("K137").Value = 3672
$ws.Range("L137").Value = 5800.0002
$ws.Range("M137").Value = -1122
$ws.Range("N137").Value = -10900.0002

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4460
$ws.Range("I63").Value = 2100
$ws.Range("J63").Value = 8000
$ws.Range("K63").Value = 2100
$ws.Range("L63").Value = 8000
$ws.Range("M63").Value = -1414
$ws.Range("N63").Value = -9372
$ws.Range("H66").Value = 4460
$ws.Range("I66").Value = 2100
$ws.Range("J66").Value = 8000
$ws.Range("K66").Value = 10500
$ws.Range("L66").Value = 40000
$ws.Range("M66").Value = -7068
$ws.Range("N66").Value = -46864
$ws.Range("H74").Value = 1271.9048
$ws.Range("I74").Value = 1051.3334
$ws.Range("J74").Value = 1566
$ws.Range("K74").Value = 1051.3334
$ws.Range("L74").Value = 1566
$ws.Range("M74").Value = -177.3334
$ws.Range("N74").Value = -3314
$ws.Range("H77").Value = 1271.9048
$ws.Range("I77").Value = 1051.3334
$ws.Range("J77").Value = 1566
$ws.Range("K77").Value = 5256.666999999999
$ws.Range("L77").Value = 7830
$ws.Range("M77").Value = -888.6669999999995
$ws.Range("N77").Value = -16566

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 46791.957
$ws.Range("I86").Value = 2453.6365
$ws.Range("K86").Value = 2453.6365
$ws.Range("M86").Value = -1330.6365
$ws.Range("H89").Value = 46791.957
$ws.Range("I89").Value = 2453.6365
$ws.Range("K89").Value = 12268.1825
$ws.Range("M89").Value = -6652.182500000001
$ws.Range("H105").Value = 4709.2856
$ws.Range("I105").Value = 2741.9
$ws.Range("K105").Value = 2741.9
$ws.Range("M105").Value = -994.9000000000001
$ws.Range("H134").Value = 8273.360000000001
$ws.Range("I134").Value = 2701.0527
$ws.Range("J134").Value = 25919
$ws.Range("K134").Value = 8103.158100000001
$ws.Range("L134").Value = 77757
$ws.Range("M134").Value = -5568.158100000001
$ws.Range("N134").Value = -82827

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1985621.5
$ws.Range("I31").Value = 2165860.2
$ws.Range("J31").Value = 2993.3333
$ws.Range("K31").Value = 2165860.2
$ws.Range("L31").Value = 2993.3333
$ws.Range("M31").Value = -2165565.2
$ws.Range("N31").Value = -3583.3333
$ws.Range("H34").Value = 1985621.5
$ws.Range("I34").Value = 2165860.2
$ws.Range("J34").Value = 2993.3333
$ws.Range("K34").Value = 2165860.2
$ws.Range("L34").Value = 2993.3333
$ws.Range("M34").Value = -2165658.2
$ws.Range("N34").Value = -3397.3333
$ws.Range("H43").Value = 30000
$ws.Range("J43").Value = 30000
$ws.Range("L43").Value = 30000
$ws.Range("N43").Value = -30368
$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 953.0625
$ws.Range("I5").Value = 883.26666
$ws.Range("K5").Value = 2649.79998
$ws.Range("M5").Value = -2537.79998
$ws.Range("H117").Value = 250955
$ws.Range("J117").Value = 333940
$ws.Range("L117").Value = 1001820
$ws.Range("N117").Value = -1008704
$ws.Range("H129").Value = 12347467
$ws.Range("I129").Value = 2005.5555
$ws.Range("J129").Value = 18520198
$ws.Range("K129").Value = 6016.666499999999
$ws.Range("L129").Value = 55560594
$ws.Range("M129").Value = -1016.666499999999
$ws.Range("N129").Value = -55570594
$ws.Range("H135").Value = 953.0625
$ws.Range("I135").Value = 883.26666
$ws.Range("K135").Value = 7949.39994
$ws.Range("M135").Value = -5414.39994

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14808247
$ws.Range("I70").Value = 29609822
$ws.Range("K70").Value = 29609822
$ws.Range("M70").Value = -29609552
$ws.Range("H73").Value = 14808247
$ws.Range("I73").Value = 29609822
$ws.Range("K73").Value = 29609822
$ws.Range("M73").Value = -29608886
$ws.Range("H102").Value = 1488.2174
$ws.Range("I102").Value = 979.125
$ws.Range("J102").Value = 2651.8572
$ws.Range("K102").Value = 979.125
$ws.Range("L102").Value = 2651.8572
$ws.Range("M102").Value = 642.875
$ws.Range("N102").Value = -5895.8572

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3032967.5
$ws.Range("I132").Value = 6668776
$ws.Range("J132").Value = 3126.7222
$ws.Range("K132").Value = 20006328
$ws.Range("L132").Value = 9380.1666
$ws.Range("M132").Value = -20003798
$ws.Range("N132").Value = -14440.1666

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 5000
$ws.Range("J29").Value = 5000
$ws.Range("L29").Value = 5000
$ws.Range("N29").Value = -5580
$ws.Range("H136").Value = 63599.938
$ws.Range("I136").Value = 101009.9
$ws.Range("J136").Value = 1250
$ws.Range("K136").Value = 303029.7
$ws.Range("L136").Value = 3750
$ws.Range("M136").Value = -300479.7
$ws.Range("N136").Value = -8850

Write-Host "All updates applied."